# "code cleanup database page"
#
# 1) Architecture slide (slide 5): nudge the picture's position slightly
#    (offset only, size is untouched).
# 2) Conclusions slide (slide 9): add a new lead-in bullet "From a project
#    to a product" above the existing first paragraph.

$p = $ppt.ActivePresentation

# --- Slide 5 ("Architecture"): reposition the screenshot picture -----------
$slide5 = $p.Slides.Item(5)
$pic = $slide5.Shapes.Item(2)   # "Content Placeholder 7" (the picture)

# Shape.Left/Top are in points; the target EMU offsets are 1393639/1409052
# (914400 EMU per inch, 12700 EMU per point). Values chosen so the
# point -> EMU round-trip lands exactly on the target offsets.
$pic.Left = 109.7354
$pic.Top = 110.94902

# --- Slide 9 ("Conclusions and future development"): add intro line -------
$slide9 = $p.Slides.Item(9)
$body = $slide9.Shapes.Item(2)  # "Content Placeholder 2" (the bullet list)

[void]$body.TextFrame.TextRange.InsertBefore("From a project to a product`r")
